# Bug fixes on parameter variation, added parameters for surface, added ACF, PSD, RAPSD
#
# The "recalcRoughL" column (H) is flipped on for every material row (it was
# left at 0.0 for all rows before; every row should now recalc the rough L),
# and the "500 nm / reverse=1 / roughdim=40" row (row 11, AlInP) is reset back
# to defaults (300 nm, reverse=0, roughdim=0) while the InGaP row (row 14)
# now gets the reverse+roughdim settings that AlInP used to carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 (AlInP): L 500 -> 300, input(reverse) 1 -> 0, roughdim 40 -> 0
$ws.Range("B5").Value = 300.0
$ws.Range("C5").Value = 0.0
$ws.Range("D5").Value = 0.0

# Row 14 (InGaP): input(reverse) 0 -> 1, roughdim 0 -> 40
$ws.Range("C9").Value = 1.0
$ws.Range("D9").Value = 40.0

# recalcRoughL (H) turned on for every data row
$ws.Range("H2").Value = 1.0
$ws.Range("H3").Value = 1.0
$ws.Range("H4").Value = 1.0
$ws.Range("H5").Value = 1.0
$ws.Range("H6").Value = 1.0
$ws.Range("H7").Value = 1.0
$ws.Range("H8").Value = 1.0
$ws.Range("H9").Value = 1.0
$ws.Range("H10").Value = 1.0
